# Update gh-pages to output generated at 456a3b4
#
# Summary of the change being applied:
#  1. On sheet "展览" (Sheet1): bump the "想去人数" (F column) value for a
#     set of events (identified by their bilibili id embedded in the H
#     column link) to reflect newly scraped counts.
#  2. On sheet "演出" (Sheet2): the event "COS STAR 次元之夜ACG主题派对3.0"
#     (2024-05-01, id=84288) was removed from the source feed. All rows
#     below it shift up by one (columns B:I only -- the leading index
#     column A keeps its original per-row numbering), and the last row
#     disappears.
#  3. On sheet "全部类型" (Sheet4, the union of all the other sheets): the
#     same event is removed, with the same up-shift behaviour, and the
#     F column bumps from (1) are re-applied (by id) at their new row
#     positions since this sheet carries its own independent copy of the
#     data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Map of bilibili "id=" -> updated "想去人数" (F column) value.
# ---------------------------------------------------------------------
$fMap = @{}
$fMap["82979"] = 1694
$fMap["83221"] = 1118
$fMap["81566"] = 1558
$fMap["85052"] = 9
$fMap["83132"] = 1475
$fMap["81962"] = 3104
$fMap["84177"] = 648
$fMap["82458"] = 1787
$fMap["83522"] = 862
$fMap["83226"] = 281
$fMap["82918"] = 1485
$fMap["84890"] = 14
$fMap["82319"] = 1232
$fMap["83462"] = 466
$fMap["84662"] = 134
$fMap["83856"] = 4843
$fMap["85020"] = 64
$fMap["83910"] = 580
$fMap["84896"] = 74
$fMap["84815"] = 165

function Get-BilibiliId($text) {
    if ($text -match "id=(\d+)") {
        return $matches[1]
    }
    return $null
}

# ---------------------------------------------------------------------
# Step 1: sheet "展览" -- update the F column (想去人数) in place by id.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$lastRow1 = $ws1.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow1; $r++) {
    $id = Get-BilibiliId $ws1.Cells.Item($r, 8).Value2
    if ($id -ne $null -and $fMap.ContainsKey($id)) {
        $ws1.Cells.Item($r, 6).Value2 = $fMap[$id]
    }
}

# ---------------------------------------------------------------------
# Step 2: sheet "演出" -- remove the "COS STAR" row (2024-05-01,
# id=84288), shifting B:I of every following row up by one, and drop
# the now-duplicated last row.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$lastRow2 = $ws2.UsedRange.Rows.Count
for ($r = 2; $r -lt $lastRow2; $r++) {
    $src = $ws2.Range("B" + ($r + 1) + ":I" + ($r + 1))
    $dst = $ws2.Range("B" + $r + ":I" + $r)
    $src.Copy($dst)
}
$ws2.Rows.Item($lastRow2).Delete()

# ---------------------------------------------------------------------
# Step 3: sheet "全部类型" -- same row removal/shift as step 2 (this
# sheet carries its own copy of every event), then re-apply the F
# column bumps from step 1 at their new row positions.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$lastRow4 = $ws4.UsedRange.Rows.Count

# Find the row holding the event that needs to be dropped (id=84288)
# rather than assuming a fixed row number.
$deleteRow = $null
for ($r = 2; $r -le $lastRow4; $r++) {
    $id = Get-BilibiliId $ws4.Cells.Item($r, 8).Value2
    if ($id -eq "84288") {
        $deleteRow = $r
        break
    }
}

if ($deleteRow -ne $null) {
    for ($r = $deleteRow; $r -lt $lastRow4; $r++) {
        $src = $ws4.Range("B" + ($r + 1) + ":I" + ($r + 1))
        $dst = $ws4.Range("B" + $r + ":I" + $r)
        $src.Copy($dst)
    }
    $ws4.Rows.Item($lastRow4).Delete()
}

$lastRow4after = $ws4.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow4after; $r++) {
    $id = Get-BilibiliId $ws4.Cells.Item($r, 8).Value2
    if ($id -ne $null -and $fMap.ContainsKey($id)) {
        $ws4.Cells.Item($r, 6).Value2 = $fMap[$id]
    }
}
